$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block updates ---
# VALOR MORA total (E11): 264000 -> 408000
$ws.Range("E11").Value = 408000
# Cant. Periodos (F13): 2 -> 3
$ws.Range("F13").Value = 3

# --- Expand the data table from 4 rows (2 workers x 2 periods) to 6 rows (2 workers x 3 periods) ---
# Insert two new rows just above the last (specially styled) data row so the
# table grows from B16:J19 to B16:J21, keeping the bottom border row intact.
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Insert()

# Copy the normal-row formatting down onto the two freshly inserted rows.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Re-populate the table, now grouped by period (2506, 2507, 2508) ---
$data = @(
  @(16, "CC", "1063495346", "JANER ENRIQUE CHAMORRO ROCHA",   "2506", 66667, 2000000),
  @(17, "CC", "1124503428", "NEIVER DE JESUS CHAMORRO ROCHA", "2506", 53333, 1600000),
  @(18, "CC", "1063495346", "JANER ENRIQUE CHAMORRO ROCHA",   "2507", 80000, 2000000),
  @(19, "CC", "1124503428", "NEIVER DE JESUS CHAMORRO ROCHA", "2507", 64000, 1600000),
  @(20, "CC", "1063495346", "JANER ENRIQUE CHAMORRO ROCHA",   "2508", 80000, 2000000),
  @(21, "CC", "1124503428", "NEIVER DE JESUS CHAMORRO ROCHA", "2508", 64000, 1600000)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
  $ws.Range("E$r").Value = $row[4]
  $ws.Range("F$r").Value = $row[5]
  $ws.Range("G$r").Value = $row[6]
}

# --- Nudge the logo image slightly to the left (cosmetic repositioning) ---
$logo = $ws.Shapes.Item(1)
$logo.Left = $logo.Left - 13.5
